$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("H2").Value = 82
$ws.Range("G3").Value = 94
$ws.Range("H3").Value = 98
$ws.Range("I3").Value = 149
$ws.Range("J3").Value = 153
$ws.Range("B6").Value = 281
$ws.Range("D6").Value = 310
$ws.Range("E6").Value = 322
$ws.Range("G6").Value = 360
$ws.Range("I6").Value = 388
$ws.Range("B7").Value = 382
$ws.Range("D7").Value = 483
$ws.Range("E7").Value = 481
$ws.Range("G7").Value = 522
$ws.Range("H7").Value = 518
$ws.Range("I7").Value = 646
$ws.Range("J7").Value = 552

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("D8").Value = 23
$ws.Range("I27").Value = 10
$ws.Range("G28").Value = 33
$ws.Range("E36").Value = 30
$ws.Range("J36").Value = 31
$ws.Range("H41").Value = 5
$ws.Range("G45").Value = 5
$ws.Range("G53").Value = 63
$ws.Range("I53").Value = 105
$ws.Range("H54").Value = 4
$ws.Range("I76").Value = 18
$ws.Range("D86").Value = 3
$ws.Range("G86").Value = 20
$ws.Range("B89").Value = 3
$ws.Range("B98").Value = 382
$ws.Range("D98").Value = 483
$ws.Range("E98").Value = 481
$ws.Range("G98").Value = 522
$ws.Range("H98").Value = 518
$ws.Range("I98").Value = 646
$ws.Range("J98").Value = 552

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("I3").Value = 2
$ws.Range("I6").Value = 18

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("D5").Value = 15
$ws.Range("D6").Value = 23

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J3").Value = 8
$ws.Range("E6").Value = 18
$ws.Range("E7").Value = 30
$ws.Range("J7").Value = 31

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("H3").Value = 3
$ws.Range("H5").Value = 5

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("G3").Value = 5
$ws.Range("D5").Value = 1
$ws.Range("D6").Value = 3
$ws.Range("G6").Value = 20

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("G6").Value = 24
$ws.Range("G7").Value = 33

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("G6").Value = 44
$ws.Range("I6").Value = 67
$ws.Range("G7").Value = 63
$ws.Range("I7").Value = 105

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("I3").Value = 2
$ws.Range("I5").Value = 10

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("B3").Value = 3
$ws.Range("B4").Value = 3

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("H2").Value = 2
$ws.Range("H5").Value = 4

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 5
